# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet, a new (blank) column was inserted
# immediately before column N ("Late"), pushing the existing
# N/O/P ("Late" / "heading" / "Outstanding") columns one place to the
# right (-> O/P/Q). The new column inherits the column-width formatting
# of the column to its left (M), as Excel does automatically on
# Insert. The workbook was left with the "Repayment schedule" tab
# active/selected, with cell R6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (moves tabSelected / activeTab onto it).
$ws.Activate()

# Insert a new blank column before column N (14th column : A=1 ... N=14).
$ws.Columns.Item(14).Insert()

# The freshly inserted column picks up the width of its left neighbour
# (column M), matching Excel's native "insert column" behaviour.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Leave the selection on R6, as recorded after the edit.
[void]$ws.Range("R6").Select()
